$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("Q2").Value = "id1-leonid.png"

# Row 3
$ws.Range("M3").Value = 4
$ws.Range("Q3").Value = "id2-elena.png"

# Row 4
$ws.Range("M4").Value = 5
$ws.Range("Q4").Value = "id3-sergey.png"

# Row 5
$ws.Range("M5").Value = 5

# Row 6
$ws.Range("M6").Value = 5
$ws.Range("Q6").Value = "id5-maksim.png"

# Row 7
$ws.Range("M7").Value = 5
$ws.Range("Q7").Value = "id6-sophia.png"

# Row 8
$ws.Range("M8").Value = 6
$ws.Range("Q8").Value = "id7-sviatik.png"

# Row 9
$ws.Range("M9").Value = 5

# Row 10
$ws.Range("M10").Value = 5

# Row 11
$ws.Range("M11").Value = 4

# Row 12
$ws.Range("M12").Value = 4
$ws.Range("Q12").Value = "id11-milana.png"

# Row 13
$ws.Range("M13").Value = 6

# Row 14
$ws.Range("M14").Value = 6

# Row 15
$ws.Range("M15").Value = 4

# Row 16
$ws.Range("M16").Value = 4

# Row 17
$ws.Range("M17").Value = 5

# Row 18
$ws.Range("M18").Value = 5

# Row 19
$ws.Range("M19").Value = 4

# Row 20
$ws.Range("M20").Value = 4

# Row 21
$ws.Range("M21").Value = 5

# Row 22
$ws.Range("M22").Value = 4

# Row 23
$ws.Range("M23").Value = 4

# Row 24
$ws.Range("M24").Value = 5

# Row 25
$ws.Range("M25").Value = 5

# Row 26
$ws.Range("M26").Value = 6

# Row 27
$ws.Range("M27").Value = 6

# Row 28
$ws.Range("M28").Value = 3

# Row 29
$ws.Range("M29").Value = 3

# Row 30
$ws.Range("M30").Value = 3
$ws.Range("Q30").Value = "id29-lilina.png"

# Row 31
$ws.Range("M31").Value = 3

# Row 32
$ws.Range("M32").Value = 2
$ws.Range("Q32").Value = "id31-sofia.png"

# Row 33
$ws.Range("M33").Value = 2
$ws.Range("Q33").Value = "id32-vasilii.png"

# Row 34
$ws.Range("M34").Value = 3
$ws.Range("Q34").Value = "id33-fedor.png"

# Row 35
$ws.Range("M35").Value = 3
$ws.Range("Q35").Value = "id34-lelia.png"

# Row 36
$ws.Range("M36").Value = 3
$ws.Range("Q36").Value = "id35-nikolay.png"

# Row 37
$ws.Range("M37").Value = 3

# Row 38
$ws.Range("M38").Value = 4

# Row 39
$ws.Range("M39").Value = 3
$ws.Range("Q39").Value = "id38-ivan.png"

# Row 40
$ws.Range("M40").Value = 4

# Row 41
$ws.Range("M41").Value = 4

# Row 42
$ws.Range("M42").Value = 4

# Row 43
$ws.Range("M43").Value = 5

# Row 44
$ws.Range("M44").Value = 4

# Row 45
$ws.Range("M45").Value = 3

# Row 46
$ws.Range("M46").Value = 4

# Row 47
$ws.Range("M47").Value = 4

# Row 48
$ws.Range("M48").Value = 5

# Row 49
$ws.Range("M49").Value = 4

# Row 50
$ws.Range("M50").Value = 5

# Row 51
$ws.Range("M51").Value = 5

# Row 52
$ws.Range("M52").Value = 1
$ws.Range("P52").Value = 'Фотография Каземира сделана  в г. Ровны в 1928 году'
$ws.Range("Q52").Value = "id51-kazimir.png"

# Row 53
$ws.Range("M53").Value = 2
$ws.Range("Q53").Value = "id52-fedor.png"

# Row 54
$ws.Range("M54").Value = 4

# Row 55
$ws.Range("M55").Value = 5

# Row 56
$ws.Range("M56").Value = 5

Write-Output "Edit complete"